$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("used")

# Remove the used entry from the names sheet (row 2: "eao5ane8") and shift
# the remaining rows up, same as Excel's native row delete.
$ws1.Rows.Item(2).Delete()

# Append the corresponding record to the "used" sheet.
$lastRow = $ws2.Cells.Item($ws2.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1
$ws2.Cells.Item($newRow, 1).Value = "eao5ane8"
$ws2.Cells.Item($newRow, 2).Value = "ChatGPT Image 2026年1月22日 23_09_10.png"
$ws2.Cells.Item($newRow, 3).Value = "2026-01-22 23:11:12"
